# Update the "想去人数" (want-to-go count) figures for the first two
# convention entries on both the "展览" and "全部类型" sheets, reflecting
# refreshed data from the generator run at 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 468
    $ws.Range("F3").Value = 3310
}
